$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2
$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 5
$ws.Range("B9").Value = 36
$ws.Range("B10").Value = 63
$ws.Range("B11").Value = 6
$ws.Range("B12").Value = 4

$ws.Range("B13").Select()
